$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions scheduled update)

# Row 2
$ws.Range("D2").Value = '64.060.63'
$ws.Range("E2").Value = '  -3.63%  '

# Row 3
$ws.Range("D3").Value = '3.163.57'
$ws.Range("E3").Value = '  -8.77%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.05'
$ws.Range("E5").Value = '  -4.23%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.91'
$ws.Range("E6").Value = '  -5.14%  '

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("E8").Value = '  -3.79%  '

# Row 9
$ws.Range("D9").Value = '3.164.07'
$ws.Range("E9").Value = '  -8.70%  '

# Row 10
$ws.Range("E10").Value = '  -7.41%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.59'
$ws.Range("E11").Value = '  -5.46%  '

# Row 12
$ws.Range("E12").Value = '  -5.90%  '

# Row 13
$ws.Range("D13").Value = '3.711.32'
$ws.Range("E13").Value = '  -8.77%  '

# Row 14
$ws.Range("E14").Value = '  +1.28%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.18'
$ws.Range("E15").Value = '  -10.34%  '

# Row 16
$ws.Range("D16").Value = '64.081.68'

# Row 17
$ws.Range("E17").Value = '  -5.87%  '

# Row 18
$ws.Range("D18").Value = '3.164.69'
$ws.Range("E18").Value = '  -8.61%  '

# Row 19
$ws.Range("E19").Value = '  -4.64%  '

# Row 20
$ws.Range("E20").Value = '  -6.68%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '351.05'
$ws.Range("E21").Value = '  -5.99%  '

# Row 22
$ws.Range("E22").Value = '  -6.78%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.53'
$ws.Range("E24").Value = '  -6.53%  '

# Row 25
$ws.Range("E25").Value = '  -6.46%  '

# Row 26
$ws.Range("E26").Value = '  -6.79%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.49'
$ws.Range("E27").Value = '  -5.05%  '

# Row 28
$ws.Range("E28").Value = '  -1.28%  '

# Row 29
$ws.Range("E29").Value = '  +0.06%  '

# Row 30
$ws.Range("E30").Value = '  -0.14%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.48'
$ws.Range("E31").Value = '  -7.73%  '

# Row 32
$ws.Range("E32").Value = '  -5.77%  '

# Row 33
$ws.Range("E33").Value = '  -7.69%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.62'
$ws.Range("E34").Value = '  -6.35%  '

# Row 35
$ws.Range("E35").Value = '  -6.08%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.44'
$ws.Range("E36").Value = '  -7.67%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '153.68'
$ws.Range("E37").Value = '  -4.70%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.814'
$ws.Range("E38").Value = '  -8.19%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '25.66'
$ws.Range("E39").Value = '  -9.03%  '

# Row 40
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.69'
$ws.Range("E40").Value = '  -6.71%  '

# Row 41
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.49'
$ws.Range("E41").Value = '  -3.63%  '

# Row 42
$ws.Range("D42").Value = '2.595.80'
$ws.Range("E42").Value = '  -6.88%  '

# Row 43
$ws.Range("E43").Value = '  -7.81%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.36'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.95'
$ws.Range("E45").Value = '  -7.84%  '

# Row 46
$ws.Range("E46").Value = '  -6.55%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.61'
$ws.Range("E47").Value = '  -6.65%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '316.96'
$ws.Range("E48").Value = '  -7.12%  '

# Row 49
$ws.Range("E49").Value = '  -8.48%  '

# Row 50
$ws.Range("E50").Value = '  -4.08%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  +0.02%  '
